$wb = $excel.ActiveWorkbook

# --- NC2: duplicate of the "NC" (No control) results sheet, appended at the end ---
$srcNC = $wb.Worksheets.Item("NC")
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcNC.Copy($null, $afterSheet)
$wsNC2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNC2.Name = "NC2"

# --- HC1: duplicate of the "HC" (Holding control) results sheet, appended at the end, ---
# --- with the refreshed numeric results from the new (longer) simulation run ---
$srcHC = $wb.Worksheets.Item("HC")
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcHC.Copy($null, $afterSheet2)
$wsHC1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHC1.Name = "HC1"

$wsHC1.Range("B2").Value = 2338.400307218047
$wsHC1.Range("C2").Value = 12390.70170756176
$wsHC1.Range("D2").Value = 132.8640879770229
$wsHC1.Range("E2").Value = 14861.96610275683
